# Edit the "API Endpoints" workbook:
#  - insert 2 new rows for the new /auth/sendVerify and /auth/verifyEmail routes
#  - insert a new "Response" column (HTTP status code) before the old Status/Implemented column
#  - rename the old Status column header to "Implemented"
#  - update several routes' data
#  - refresh conditional formatting to point at the relocated column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Structural changes: make room for the new rows/column first.
# ---------------------------------------------------------------------------
# Two new rows for /auth/sendVerify (POST) and /auth/verifyEmail (GET),
# inserted right after the existing /auth/refresh row (old row 9).
$ws.Rows.Item(10).Resize(2).Insert()

# A new column for the HTTP response code, inserted before the old
# Status/Implemented column (old column E, now shifted to F).
$ws.Columns.Item(5).Insert()

# ---------------------------------------------------------------------------
# 2) Fill in the brand-new content in the same order the author originally
#    typed it, so freshly-created shared strings land in the right slots.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "/auth/sendVerify"
$ws.Range("C10").Value = "POST"
$ws.Range("D10").Value = "Send verification email"

$ws.Range("B11").Value = "/auth/verifyEmail"
$ws.Range("C11").Value = "GET"
$ws.Range("D11").Value = "Verify email via tokenized link"

# The admin "list a user's borrowed songs" route moved from an {id} style
# param to an {email} style param.
$ws.Range("B18").Value = "/admin/users/{email}/lend"

# DELETE /user has no meaningful response code yet - entered as literal text
# "-" (quote-prefixed so Excel keeps it as text) and right aligned.
$ws.Range("E5").Value = "'-"
$ws.Range("E5").HorizontalAlignment = -4152

# New column headers.
$ws.Range("F1").Value = "Implemented"
$ws.Range("E1").Value = "Response"

# ---------------------------------------------------------------------------
# 3) Fill in the rest of the new "Response" column (HTTP status codes).
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = 201
$ws.Range("E3").Value = 200
$ws.Range("E4").Value = 200
$ws.Range("E6").Value = 201
$ws.Range("E7").Value = 200
$ws.Range("E8").Value = 201
$ws.Range("E9").Value = 201
$ws.Range("E10").Value = 200
$ws.Range("E11").Value = 302
$ws.Range("E12").Value = 200
$ws.Range("E13").Value = 200
$ws.Range("E14").Value = 200
$ws.Range("E15").Value = 200
$ws.Range("E16").Value = 200
$ws.Range("E17").Value = 200
$ws.Range("E18").Value = 200
$ws.Range("E19").Value = 200
$ws.Range("E20").Value = 200

# ---------------------------------------------------------------------------
# 4) A handful of rows also had their "Implemented" status upgraded.
# ---------------------------------------------------------------------------
$ws.Range("F4").Value = "Fully Implemented"
$ws.Range("F7").Value = "Fully Implemented"
$ws.Range("F8").Value = "Fully Implemented"

# ---------------------------------------------------------------------------
# 5) Column widths adjust slightly to fit the new layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 24.85546875
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 17.42578125

# ---------------------------------------------------------------------------
# 6) Conditional formatting needs to follow the relocated "Implemented"
#    column (now F instead of E).
# ---------------------------------------------------------------------------
$condRange = $ws.Range("E2:E20")
$fcs = $condRange.FormatConditions
$rule1 = $fcs.Item(1)
$rule1.Formula1 = '=F2="Fully Implemented"'
$rule1.ModifyAppliesToRange($ws.Range("F2:F20"))

$rule2 = $ws.Range("F2:F20").FormatConditions.Item(2)
$rule2.Formula1 = '=F2="POC / MVP"'
$rule2.ModifyAppliesToRange($ws.Range("F2:F20"))

$rule3 = $ws.Range("F2:F20").FormatConditions.Item(3)
$rule3.Formula1 = '=F2="Not Implemented"'
$rule3.ModifyAppliesToRange($ws.Range("F2:F20"))
